$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 287.92307
$ws.Range("I28").Value = 246.18182
$ws.Range("K28").Value = 246.18182
$ws.Range("M28").Value = 238.81818

$ws.Range("H86").Value = 57975220
$ws.Range("I86").Value = 83337490
$ws.Range("J86").Value = 30307294
$ws.Range("K86").Value = 83337490
$ws.Range("L86").Value = 30307294
$ws.Range("M86").Value = -83336367
$ws.Range("N86").Value = -30309540

$ws.Range("H89").Value = 57975220
$ws.Range("I89").Value = 83337490
$ws.Range("J89").Value = 30307294
$ws.Range("K89").Value = 416687450
$ws.Range("L89").Value = 151536470
$ws.Range("M89").Value = -416681834
$ws.Range("N89").Value = -151547702

$ws.Range("H94").Value = 1337.8182
$ws.Range("I94").Value = 1337.8182
$ws.Range("K94").Value = 1337.8182
$ws.Range("M94").Value = -886.8181999999999

$ws.Range("H106").Value = 52634160
$ws.Range("I106").Value = 55558030
$ws.Range("J106").Value = 4444
$ws.Range("K106").Value = 55558030
$ws.Range("L106").Value = 4444
$ws.Range("M106").Value = -55557399
$ws.Range("N106").Value = -5706

$ws.Range("H132").Value = 1006.59375
$ws.Range("I132").Value = 974.6786
$ws.Range("K132").Value = 2924.0358
$ws.Range("M132").Value = -394.0357999999997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").Value = $null

$ws.Range("H110").Value = 2511.16
$ws.Range("I110").Value = 1454.7222
$ws.Range("J110").Value = 5227.7144
$ws.Range("K110").Value = 1454.7222
$ws.Range("L110").Value = 5227.7144
$ws.Range("M110").Value = 590.2778000000001
$ws.Range("N110").Value = -9317.7144

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 76924130
$ws.Range("I80").Value = 125001420
$ws.Range("K80").Value = 125001420
$ws.Range("M80").Value = -125000422

$ws.Range("H83").Value = 76924130
$ws.Range("I83").Value = 125001420
$ws.Range("K83").Value = 625007100
$ws.Range("M83").Value = -625002108

$ws.Range("H99").Value = 10505.375
$ws.Range("I99").Value = 11577.571
$ws.Range("K99").Value = 11577.571
$ws.Range("M99").Value = -10079.571

$ws.Range("H107").Value = 12593.429
$ws.Range("I107").Value = 14317.958
$ws.Range("K107").Value = 14317.958
$ws.Range("M107").Value = -12397.958

$ws.Range("H134").Value = 25716132
$ws.Range("I134").Value = 1867.0476
$ws.Range("K134").Value = 5601.142800000001
$ws.Range("M134").Value = -3066.142800000001

$ws.Range("H135").Value = 74198
$ws.Range("J135").Value = 74198
$ws.Range("L135").Value = 74198
$ws.Range("N135").Value = -84338

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2932.3823
$ws.Range("I31").Value = 2651.6775
$ws.Range("K31").Value = 2651.6775
$ws.Range("M31").Value = -2356.6775

$ws.Range("H34").Value = 2932.3823
$ws.Range("I34").Value = 2651.6775
$ws.Range("K34").Value = 2651.6775
$ws.Range("M34").Value = -2449.6775

$ws.Range("H106").Value = 56415.25
$ws.Range("J106").Value = 56415.25
$ws.Range("L106").Value = 56415.25
$ws.Range("N106").Value = -58939.25

$ws.Range("H122").Value = 4423.5
$ws.Range("I122").Value = 6300
$ws.Range("K122").Value = 18900
$ws.Range("M122").Value = -16450

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = $null

$ws.Range("H132").Value = 1866.2858
$ws.Range("I132").Value = 1504.4
$ws.Range("K132").Value = 4513.200000000001
$ws.Range("M132").Value = -1983.200000000001

$ws.Range("H134").Value = 2068.087
$ws.Range("I134").Value = 1848.2941
$ws.Range("J134").Value = 2690.8333
$ws.Range("K134").Value = 5544.8823
$ws.Range("L134").Value = 8072.499899999999
$ws.Range("M134").Value = -3009.8823
$ws.Range("N134").Value = -13142.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 36050.867
$ws.Range("I44").Value = 397.5
$ws.Range("J44").Value = 41536
$ws.Range("K44").Value = 1192.5
$ws.Range("L44").Value = 124608
$ws.Range("M44").Value = -794.5
$ws.Range("N44").Value = -125404

$ws.Range("H80").Value = 1899.25
$ws.Range("J80").Value = 1549
$ws.Range("L80").Value = 4647
$ws.Range("N80").Value = -6519

$ws.Range("H83").Value = 1899.25
$ws.Range("J83").Value = 1549
$ws.Range("L83").Value = 13941
$ws.Range("N83").Value = -23301

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 55573644
$ws.Range("J80").Value = 166700670
$ws.Range("L80").Value = 166700670
$ws.Range("N80").Value = -166702666

$ws.Range("H83").Value = 55573644
$ws.Range("J83").Value = 166700670
$ws.Range("L83").Value = 833503350
$ws.Range("N83").Value = -833513334

$ws.Range("H102").Value = 41668044
$ws.Range("I102").Value = 50001256
$ws.Range("K102").Value = 50001256
$ws.Range("M102").Value = -49999634

$ws.Range("I107").Value = 250392
$ws.Range("J107").Value = 1355.4
$ws.Range("K107").Value = 250392
$ws.Range("L107").Value = 1355.4
$ws.Range("M107").Value = -248472
$ws.Range("N107").Value = -5195.4

$ws.Range("H122").Value = 2824.6875
$ws.Range("I122").Value = 2092.3845
$ws.Range("J122").Value = 5998
$ws.Range("K122").Value = 6277.1535
$ws.Range("L122").Value = 17994
$ws.Range("M122").Value = -3827.1535
$ws.Range("N122").Value = -22894

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2181.3447
$ws.Range("I40").Value = 2282.0454
$ws.Range("J40").Value = 1864.8572
$ws.Range("K40").Value = 2282.0454
$ws.Range("L40").Value = 1864.8572
$ws.Range("M40").Value = -2146.0454
$ws.Range("N40").Value = -2136.8572

$ws.Range("H68").Value = 7928
$ws.Range("I68").Value = 10299.2
$ws.Range("K68").Value = 10299.2
$ws.Range("M68").Value = -9550.2

$ws.Range("H71").Value = 7928
$ws.Range("I71").Value = 10299.2
$ws.Range("K71").Value = 51496
$ws.Range("M71").Value = -47752

$ws.Range("H82").Value = 1564.7
$ws.Range("I82").Value = 1807.125
$ws.Range("J82").Value = 595
$ws.Range("K82").Value = 1807.125
$ws.Range("L82").Value = 595
$ws.Range("M82").Value = -1446.125
$ws.Range("N82").Value = -1317

$ws.Range("H85").Value = 1564.7
$ws.Range("I85").Value = 1807.125
$ws.Range("J85").Value = 595
$ws.Range("K85").Value = 1807.125
$ws.Range("L85").Value = 595
$ws.Range("M85").Value = -559.125
$ws.Range("N85").Value = -3091

$ws.Range("H101").Value = 21677.285
$ws.Range("J101").Value = 21677.285
$ws.Range("L101").Value = 21677.285
$ws.Range("N101").Value = -28167.285

$ws.Range("H103").Value = 23128.666
$ws.Range("J103").Value = 23128.666
$ws.Range("L103").Value = 23128.666
$ws.Range("N103").Value = -25472.666

$ws.Range("H122").Value = 3054.125
$ws.Range("I122").Value = 2782
$ws.Range("K122").Value = 8346
$ws.Range("M122").Value = -5896

$ws.Range("H132").Value = 5710.364
$ws.Range("I132").Value = 5227.375
$ws.Range("K132").Value = 15682.125
$ws.Range("M132").Value = -13152.125

$ws.Range("H136").Value = 54898.21
$ws.Range("I136").Value = 61062.707
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 183188.121
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -180638.121
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 999999
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 999999
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 1999998
$ws.Range("M81").Value = $null
$ws.Range("N81").Value = -2002120

$ws.Range("H84").Value = 999999
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 999999
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 9999990
$ws.Range("M84").Value = $null
$ws.Range("N84").Value = -10010598

$ws.Range("H96").Value = 3446.3333
$ws.Range("I96").Value = 3355.111
$ws.Range("K96").Value = 3355.111
$ws.Range("M96").Value = -1982.111

$ws.Range("H101").Value = 29638.75
$ws.Range("J101").Value = 29638.75
$ws.Range("L101").Value = 29638.75
$ws.Range("N101").Value = -36128.75

$ws.Range("H113").Value = 817
$ws.Range("I113").Value = 774.75
$ws.Range("J113").Value = 901.5
$ws.Range("K113").Value = 2324.25
$ws.Range("L113").Value = 2704.5
$ws.Range("M113").Value = -154.25
$ws.Range("N113").Value = -7044.5

$ws.Range("H126").Value = 7815760
$ws.Range("I126").Value = 10003031
$ws.Range("K126").Value = 30009093
$ws.Range("M126").Value = -30006623

$ws.Range("H132").Value = 2872.6365
$ws.Range("I132").Value = 1967
$ws.Range("K132").Value = 5901
$ws.Range("M132").Value = -3371
